# Updated task estimation for WalletBalance user story.
# Adds a new "GetBalance" user story block (rows 39-43) to Sheet1,
# mirroring the structure of the other user-story sections already
# present in the estimation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user story header (row 39, column A)
$ws.Range("A39").Value = "Kao korisnik potrebno je da mogu da proverim stanje na svom novcaniku u svakom trenutku"

# Tasks + estimated time (minutes) for the new user story (rows 40-43)
$ws.Range("B40").Value = "Dodavanje GetBalance metode na WalletService"
$ws.Range("C40").Value = 5

$ws.Range("B41").Value = "Implementacija testova za GetBalance"
$ws.Range("C41").Value = 15

$ws.Range("B42").Value = "Dodavanje rute za proveru stanja Walleta u WalletController"
$ws.Range("C42").Value = 10

$ws.Range("B43").Value = "Dodavanje stranice za proveru stanja Walleta u MVC aplikaciju"
$ws.Range("C43").Value = 20

# Leave the selection where the author left it after entering the new rows
$ws.Range("A45").Select()
